$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capitalize "blue" -> "Blue" in the relevant cells' text content.

$ws.Range("F7").Value = "color: `"Blue`"`nside_1 : 5`nside_2 : 5`nside_3 : 7"

$ws.Range("F9").Value = "color : `"Blue`"`nside_1 : `"5`"`nside_2 : 5`nside_3 : 7"

$ws.Range("F10").Value = "color : `"Blue`"`nside_1 : 5`nside_2 : `"5`"`nside_3 : 7"

$ws.Range("F11").Value = "color : `"Blue`"`nside_1 : 5`nside_2 : 5`nside_3 : `"7`""

$ws.Range("E12").Value = "Triangle(`"Blue`", 5, 5, 7)"
$ws.Range("E13").Value = "Triangle(`"Blue`", 5, 5, 7)"
$ws.Range("E14").Value = "Triangle(`"Blue`", 5, 5, 7)"

$ws.Range("G12").Value = "The shape color is Blue.`nThis triangle has three sides with lengths of 5, 5, 7 centimeters."

# Update the view: move selection from G13 to G12.
$ws.Range("G12").Select()
